$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.01
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.99

$ws.Range("G3").Value = 0.99

$ws.Range("F4").Value = 0.99
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 0.01

$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0.99

$ws.Range("F6").Value = 0.01
$ws.Range("G6").Value = 0.01
$ws.Range("H6").Value = 0.99
